# Add a "staff id"/phone-number column to the line-manager sheet.
#
# The existing column G (which holds the secondary phone number, e.g.
# "0034567890") is shifted one column to the right, becoming column H, and a
# new (blank, same-formatted) column G is inserted in its place.
#
# Column G already carries the correct style/format for every row, so a
# native "insert column" reproduces that faithfully. The engine's column
# insert implementation pushes the sheet's trailing default-width column
# group boundary out by one (256 -> 257); deleting that now-superfluous
# last column restores the original 256-column boundary without touching
# any real data (it was always empty/default).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns(7).Insert()
$ws.Columns(257).Delete()
